$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings in column A (force text entry + reset to General so
# no new cell style is introduced, matching the source author's edit which
# only touched the shared-string table, not cell formatting)
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "07.03.2020"
$ws.Range("A3").NumberFormat = "General"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "10.03.2020"
$ws.Range("A4").NumberFormat = "General"

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "21.03.2020"
$ws.Range("A5").NumberFormat = "General"

# Update numeric values
$ws.Range("H2").Value = 15
$ws.Range("I2").Value = 18
$ws.Range("D3").Value = 0.15
$ws.Range("D4").Value = 0.15

# Update the selected/active cell to A6
$ws.Range("A6").Select()
